# Biomass price table: add Low/High sub-headers for Straw, Wood Chips and
# Wood Pellets in row 2 (the triplet of Low/Med/High columns under each
# fuel type previously only had the fuel name + Euro/GJ unit spelled out,
# leaving the Low/High sub-header cells blank). Also drop the stray
# "applyNumberFormat" cell style that had been applied across the whole
# table so the data reverts to the workbook's default (General) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sub-header labels for row 2 (Low/Med/High triplets live in row 3;
# the fuel-name / unit cells already in row 2 are left untouched).
$ws.Range("B2").Value = "Straw Low"
$ws.Range("D2").Value = "Straw High"
$ws.Range("E2").Value = "Wood Chips Low"
$ws.Range("G2").Value = "Wood Chips High"
$ws.Range("H2").Value = "Wood Pellets Low"
$ws.Range("J2").Value = "Wood Pellets High"

# Remove the extra number-format style from the whole table so every cell
# falls back to the workbook's default style.
$ws.Range("A1:J30").ClearFormats()

# Match the saved selection left by the editing session.
$ws.Range("J3").Select()
